$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for new columns I and J (copy H1's formatting, then set text)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for rows 2 through 67 (I = I0 column, J = IF column)
$data = @{
    2 = @(8, 8)
    3 = @(6, 6)
    4 = @(6, 7)
    5 = @(7, 8)
    6 = @(5, 7)
    7 = @(6, 7)
    8 = @(7, 7)
    9 = @(7, 7)
    10 = @(6, 6)
    11 = @(6, 6)
    12 = @(7, 8)
    13 = @(7, 7)
    14 = @(5, 6)
    15 = @(6, 6)
    16 = @(8, 8)
    17 = @(9, 9)
    18 = @(7, 8)
    19 = @(6, 6)
    20 = @(8, 8)
    21 = @(9, 9)
    22 = @(7, 7)
    23 = @(7, 7)
    24 = @(8, 8)
    25 = @(8, 8)
    26 = @(8, 8)
    27 = @(7, 8)
    28 = @(6, 6)
    29 = @(8, 8)
    30 = @(7, 7)
    31 = @(8, 8)
    32 = @(7, 7)
    33 = @(6, 7)
    34 = @(7, 7)
    35 = @(5, 6)
    36 = @(6, 7)
    37 = @(7, 8)
    38 = @(6, 6)
    39 = @(8, 8)
    40 = @(9, 9)
    41 = @(6, 6)
    42 = @(4, 6)
    43 = @(6, 6)
    44 = @(6, 7)
    45 = @(6, 7)
    46 = @(3, 4)
    47 = @(8, 9)
    48 = @(9, 9)
    49 = @(10, 10)
    50 = @(9, 9)
    51 = @(7, 7)
    52 = @(7, 7)
    53 = @(5, 6)
    54 = @(6, 7)
    55 = @(5, 5)
    56 = @(5, 6)
    57 = @(8, 8)
    58 = @(6, 6)
    59 = @(7, 7)
    60 = @(7, 7)
    61 = @(5, 5)
    62 = @(5, 5)
    63 = @(6, 6)
    64 = @(7, 7)
    65 = @(6, 6)
    66 = @(9, 9)
    67 = @(7, 7)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
